$wb = $excel.ActiveWorkbook

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# --- Restricciones_del_lider ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2").Value = "0.8 - x"
Set-TextValue $ws2 "B2" "-1.8"
Set-TextValue $ws2 "D2" "0.74"
$ws2.Range("A3").Value = "-0.8 + x"
Set-TextValue $ws2 "B3" "-0.19999999999999996"
Set-TextValue $ws2 "D3" "0.96"

# --- Restricciones_del_follower ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2").Value = "-5.664728682170543 + 3.062015503875969y"
Set-TextValue $ws3 "B2" "4.664728682170543"
Set-TextValue $ws3 "D2" "0.76"
Set-TextValue $ws3 "E2" "0.8"
Set-TextValue $ws3 "F2" "7.9"
$ws3.Range("A3").Value = "2.9230000000000005 - 1.58y"
Set-TextValue $ws3 "B3" "-3.9230000000000005"
Set-TextValue $ws3 "D3" "0.2"
Set-TextValue $ws3 "E3" "9.5"
Set-TextValue $ws3 "F3" "8.0"

# --- Punto_modificado ---
$ws4 = $wb.Worksheets.Item(4)
Set-TextValue $ws4 "A2" "0.8"
Set-TextValue $ws4 "B2" "1.85"

# --- Vector_bf ---
$ws5 = $wb.Worksheets.Item(5)
Set-TextValue $ws5 "A2" "-11.714381782945736"

# --- Vector_BF ---
$ws6 = $wb.Worksheets.Item(6)
Set-TextValue $ws6 "A2" "8.476"
Set-TextValue $ws6 "A3" "-74.09191240310076"

# --- Vector_Alpha ---
$ws7 = $wb.Worksheets.Item(7)
$ws7.Range("A2").Value = 2.58
